$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.172.82'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.43%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.829.86'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.17%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.10%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6164'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.29%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07344'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.58%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2903'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.74%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.18'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.42%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07638'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.48%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.828.65'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.23%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.977'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.49%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6706'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.21%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.37'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.37%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008982'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.09%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.842'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.12%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.158.73'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.40%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.080.23'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.12%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '236.13'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.22%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.48'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.49%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.06%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.358'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.26%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.17%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.66'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.31%  '

# Row 26
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.517'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.37%  '

# Row 27
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1386'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.69%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.64'
$ws.Range('D28').ClearFormats()

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.488'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.53%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05867'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.58%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.078'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.85%  '

# Row 32
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.220'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.09%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.083'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.56%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.856'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.13%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.135'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.33%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7230'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.20%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.616'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.48%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.863'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.43%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.225.56'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.48%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01758'
$ws.Range('D40').ClearFormats()

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.192'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.98%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9022'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.60%  '

# Row 43
$ws.Range('E43').Value = '  +0.09%  '

# Row 44
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.998.57'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.97%  '

# Row 45
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.91'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.09%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.68'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.31%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5039'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.78%  '

# Row 48
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4042'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.67%  '

# Row 49
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000117'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.26%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.146'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.45%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1144'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.72%  '
